$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows at 13:14 for a "Placeholder Name" lookup table ---
# (This pushes the existing "Science" section and the trailing blank rows down by two.)
$ws.Rows("13:14").Insert()
$ws.Rows("13:14").RowHeight = 18.75

$ws.Range("B13").Value = "Placeholder Name"
$ws.Range("C13").Value = """A"""
$ws.Range("D13").Value = """E"""

$ws.Range("B14").Value = "Placeholder Name"
$ws.Range("C14").Value = """A"""
$ws.Range("D14").Value = """F"""

# --- New "Board / IP / Netmask / Gateway" reference table in columns J:M ---
# Give the header row (row 2) the same bold/large style already used by the
# other column headers (E2:G2).
$ws.Range("E2").Copy()
$ws.Range("J2:M2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("J2").Value = "Board "
$ws.Range("K2").Value = "IP"
$ws.Range("L2").Value = "Netmask"
$ws.Range("M2").Value = "Gateway"

$ws.Range("J3").Value = "Master"
$ws.Range("K3").Value = "192.168.1.11"
$ws.Range("L3").Value = "255.255.255.0"
$ws.Range("M3").Value = "192.168.1.1"

$ws.Range("J4").Value = "Drivetrain"
$ws.Range("K4").Value = "192.168.1.12"
$ws.Range("L4").Value = "255.255.255.0"
$ws.Range("M4").Value = "192.168.1.1"

$ws.Range("J5").Value = "Science"
$ws.Range("K5").Value = "192.168.1.13"
$ws.Range("L5").Value = "255.255.255.0"
$ws.Range("M5").Value = "192.168.1.1"

$ws.Range("J6").Value = "Arm"
$ws.Range("K6").Value = "192.168.1.14"
$ws.Range("L6").Value = "255.255.255.0"
$ws.Range("M6").Value = "192.168.1.1"

$ws.Range("J7").Value = "Jeffs board"
$ws.Range("K7").Value = "Same as ma"

# --- Column width tweaks to fit the new table ---
$ws.Columns("G").ColumnWidth = 28.666666666666668
$ws.Columns("J").ColumnWidth = 23.166666666666668
$ws.Columns("K").ColumnWidth = 32.833333333333336
$ws.Columns("L").ColumnWidth = 26
$ws.Columns("M").ColumnWidth = 20

# --- View tweaks: zoom to 90% and move the active selection to G2 ---
$excel.ActiveWindow.Zoom = 90
[void]$ws.Range("G2").Select()
